$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The task "Validacion de cuit para mostrar mensaje correcto" (row 28) is
# removed entirely; all rows below it shift up by one.
$ws.Rows.Item(28).Delete()

# "Borrar reportes del disco" (now row 32, after the shift) is completed:
# reports are shown without being saved to disk, so assign it to Agustina
# and mark it done (100%).
$ws.Range("B32").Style = "Normal"
$ws.Range("B32").Value = "Agustina"
$ws.Range("C32").Value = 1
$ws.Range("C32").NumberFormat = "0%"

# Update the active selection to reflect where the user ended up working.
$ws.Range("C33").Select()
